$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember where the old last row ("Olun Elixir") lives before we shift anything.
$oldLastRow = $ws.UsedRange.Rows.Count

# Insert 4 new rows at the top (pushes all existing data down by 4 rows)
$ws.Range("A1:A4").EntireRow.Insert()

# New IDs are large integers (19 digits) that must stay text, otherwise
# double-precision rounding would corrupt them. Format the cells as Text
# before assigning so Excel doesn't coerce them to numbers.
$ws.Range("A1:A4").NumberFormat = "@"

# Fill in the new rows with the newly added resource links
$ws.Range("A1").Value = "1458690439791251467"
$ws.Range("B1").Value = "bdo-resources"
$ws.Range("C1").Value = "pvp dmg calc "

$ws.Range("A2").Value = "1421285428396425277"
$ws.Range("B2").Value = "bdo-resources"
$ws.Range("C2").Value = " cross comm guide"

$ws.Range("A3").Value = "1406978880136286219"
$ws.Range("B3").Value = "bdo-resources"
$ws.Range("C3").Value = "New armour calc"

$ws.Range("A4").Value = "1406857078886109184"
$ws.Range("B4").Value = "bdo-resources"
$ws.Range("C4").Value = "Edania cheat sheet"

# Remove the now-obsolete last row (" Olun Elixir"), which shifted down by 4 rows
# (from row 35 to row 39) because of the insert above.
$ws.Range("A" + ($oldLastRow + 4)).EntireRow.Delete()
